$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-edit) row order/content for rows 2..9, columns A-E
# (title, timestamp, historical distance, time bucket, uri, uri-link-address, uri-link-subaddress)
# The CNN row's uri has a "#fbid=..." fragment: the cell text keeps the full
# URL, but the hyperlink itself stores the fragment as a separate SubAddress
# (matches how Excel splits Target/location for such links).
$rows = @(
    @{ A = "Stampede at German Love Parade festival kills 19"; B = "2010-07-25T05:43:15UTC"; C = 1; D = "day_1"; E = "https://www.bbc.co.uk/news/world-europe-10751899"; LinkAddr = "https://www.bbc.co.uk/news/world-europe-10751899"; LinkSub = "" },
    @{ A = "Crowd Disasters as Systemic Failures: Analysis of the Love Parade Disaster"; B = "2012-06-06T00:00:00UTC"; C = 683; D = "day_31_beyond"; E = "https://ui.adsabs.harvard.edu/abs/2012arXiv1206.5856H"; LinkAddr = "https://ui.adsabs.harvard.edu/abs/2012arXiv1206.5856H"; LinkSub = "" },
    @{ A = "Love Parade tragedy: 'I'll never forget the sight of all those twisted"; B = "2010-07-25T00:00:00UTC"; C = 1; D = "day_1"; E = "https://www.independent.co.uk/news/world/europe/love-parade-tragedy-ill-never-forget-the-sight-of-all-those-twisted-bodies-2035410.html"; LinkAddr = "https://www.independent.co.uk/news/world/europe/love-parade-tragedy-ill-never-forget-the-sight-of-all-those-twisted-bodies-2035410.html"; LinkSub = "" },
    @{ A = "Australian killed in German music festival stampede"; B = "2010-07-25T18:14:00UTC"; C = 1; D = "day_1"; E = "http://www.abc.net.au/news/stories/2010/07/25/2963606.htm"; LinkAddr = "http://www.abc.net.au/news/stories/2010/07/25/2963606.htm"; LinkSub = "" },
    @{ A = "The Love Parade: European Techno, The EDM Festival, and The Tragedy in Duisburg"; B = "1-01-01T00:00:00UTC"; C = "unknown"; D = "unknown"; E = "https://www.academia.edu/12893380/The_Love_Parade_European_Techno_The_EDM_Festival_and_The_Tragedy_in_Duisburg"; LinkAddr = "https://www.academia.edu/12893380/The_Love_Parade_European_Techno_The_EDM_Festival_and_The_Tragedy_in_Duisburg"; LinkSub = "" },
    @{ A = "Nineteen Dead In German Festival Horror"; B = "1-01-01T00:00:00UTC"; C = "unknown"; D = "unknown"; E = "https://web.archive.org/web/20100730040807/http://uk.news.yahoo.com/5/20100725/twl-nineteen-dead-in-german-festival-hor-3fd0ae9.html"; LinkAddr = "https://web.archive.org/web/20100730040807/http://uk.news.yahoo.com/5/20100725/twl-nineteen-dead-in-german-festival-hor-3fd0ae9.html"; LinkSub = "" },
    @{ A = "Loveparade: RUHR.2010"; B = "1-01-01T00:00:00UTC"; C = "unknown"; D = "unknown"; E = "https://web.archive.org/web/20100730030559/http://www.essen-fuer-das-ruhrgebiet.ruhr2010.de/en/programme/a-time-to-celebrate/loveparade.html"; LinkAddr = "https://web.archive.org/web/20100730030559/http://www.essen-fuer-das-ruhrgebiet.ruhr2010.de/en/programme/a-time-to-celebrate/loveparade.html"; LinkSub = "" },
    @{ A = "18 killed amid panic at Germany's 'Love Parade'"; B = "1-01-01T00:00:00UTC"; C = "unknown"; D = "unknown"; E = "http://edition.cnn.com/2010/WORLD/europe/07/24/germany.panic.deaths/index.html?hpt=T2#fbid=UFNfp6AHpEJ"; LinkAddr = "http://edition.cnn.com/2010/WORLD/europe/07/24/germany.panic.deaths/index.html?hpt=T2"; LinkSub = "fbid=UFNfp6AHpEJ" }
)

# Preserve the existing "uri" column hyperlink style so re-adding the
# hyperlinks below doesn't leave the cells on a new (duplicate) style.
$uriStyle = $ws.Range("E2").Style

# Remove the existing hyperlinks on the uri column before rewriting values/links
$ws.Range("E2:E9").Hyperlinks.Delete()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    if ($row.LinkSub -ne "") {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), $row.LinkAddr, $row.LinkSub)
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), $row.LinkAddr)
    }
    $ws.Cells.Item($r, 5).Style = $uriStyle
    $r = $r + 1
}
